# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly refreshed counts from the source (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row => new F value
$sheet1Changes = @{
    5  = 298
    7  = 129
    13 = 2569
    14 = 94
    15 = 33
    17 = 23
    19 = 545
    20 = 618
    22 = 98
    23 = 50
    24 = 20
    26 = 2185
    27 = 4758
    32 = 249
    33 = 2150
    37 = 43
    41 = 741
    45 = 438
}

# Sheet "全部类型" (All types) - row => new F value
$sheet4Changes = @{
    5  = 298
    7  = 129
    13 = 2569
    14 = 94
    15 = 33
    18 = 23
    20 = 545
    21 = 618
    23 = 98
    24 = 50
    25 = 20
    27 = 2185
    28 = 4758
    33 = 249
    34 = 2150
    38 = 43
    42 = 741
    46 = 438
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Changes.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Changes[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Changes.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Changes[$row]
}
